$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function New-PkgFragment($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $xmlNs + '><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) "If using an Xbox One devkit..." paragraph: split the run, wrapping
#    "Gaming.Xbox.XboxOne.x" in proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$idx1 = Get-ParaIndexByText $d "*Xbox One devkit, set the active solution platform to*"
$p1 = $d.Paragraphs($idx1)
$pAttrs1 = 'w14:paraId="142064C3" w14:textId="77777777" w:rsidR="000C47E4" w:rsidRDefault="000C47E4" w:rsidP="000C47E4"'
$body1 = '<w:r><w:t xml:space="preserve">If using an Xbox One devkit, set the active solution platform to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Gaming.Xbox.XboxOne.x</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>64.</w:t></w:r>'
$frag1 = New-PkgFragment ('<w:p ' + $pAttrs1 + '>' + $body1 + '</w:p>')
$p1.Range.InsertXML($frag1)

# ---------------------------------------------------------------------------
# 2) "If using Project Scarlett..." paragraph: reworded to reference the
#    Xbox Series X|S devkit and split into several runs, wrapping
#    "Gaming.Xbox.Scarlett.x" in proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$idx2 = Get-ParaIndexByText $d "*Project Scarlett, set the active solution platform to*"
$p2 = $d.Paragraphs($idx2)
$pAttrs2 = 'w14:paraId="080B192B" w14:textId="77777777" w:rsidR="000C47E4" w:rsidRDefault="000C47E4" w:rsidP="000C47E4"'
$body2 = '<w:r><w:t xml:space="preserve">If using </w:t></w:r><w:r><w:t>an Xbox Series X|S devkit</w:t></w:r><w:r><w:t xml:space="preserve">, set the active solution platform to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Gaming.Xbox.Scarlett.x</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>64.</w:t></w:r>'
$frag2 = New-PkgFragment ('<w:p ' + $pAttrs2 + '>' + $body2 + '</w:p>')
$p2.Range.InsertXML($frag2)

# ---------------------------------------------------------------------------
# 3) XUID paragraph: split the trailing run, wrapping "who's" in proofErr
#    gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$idx3 = Get-ParaIndexByText $d "*including a XUID to indicate*"
$p3 = $d.Paragraphs($idx3)
$pAttrs3 = 'w14:paraId="15330076" w14:textId="53EC0ACF" w:rsidR="00A72B3C" w:rsidRDefault="00275A29" w:rsidP="00275A29"'
$pPr3 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="24"/></w:numPr></w:pPr>'
$body3 = '<w:r><w:t xml:space="preserve">While the implementation allows for both Global and Social leaderboards to be queried, the only real difference between them is setting an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>enum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> value passed into the query and including a XUID to indicate </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>who' + [char]0x2019 + 's</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> friends should be included in the result list.  </w:t></w:r>'
$frag3 = New-PkgFragment ('<w:p ' + $pAttrs3 + '>' + $pPr3 + $body3 + '</w:p>')
$p3.Range.InsertXML($frag3)

"done"
